$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks cleanly (collection-level delete avoids the
# duplicate-relationship bug that per-item .Delete()/.Address edits trigger)
$ws.Hyperlinks.Delete()

# Row 2: Shri Sandip Pradhan takes charge as Whole Time Member, SEBI
$ws.Range("A2").Value = "SEBI"
$ws.Range("B2").Value = "Press Release"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2025"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "December"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2025-12-05"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "Shri Sandip Pradhan takes charge as Whole Time Member, SEBI"
$ws.Range("G2").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765190360681.pdf"
$ws.Range("H2").Value = "1765190360681.pdf"
$ws.Range("I2").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Press Release/2025/December/1765190360681.pdf"

# Row 3: India: Financial Sector Assessment Program, 2024
$ws.Range("A3").Value = "SEBI"
$ws.Range("B3").Value = "Press Release"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2025"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "December"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2025-12-03"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "India: Financial Sector Assessment Program, 2024"
$ws.Range("G3").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1764758072366.pdf"
$ws.Range("H3").Value = "1764758072366.pdf"
$ws.Range("I3").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Press Release/2025/December/1764758072366.pdf"

# Row 4: Niveshak Shivir to be held on December 06, 2025 in the ci...
$ws.Range("A4").Value = "SEBI"
$ws.Range("B4").Value = "Press Release"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2025"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "December"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2025-12-02"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "Niveshak Shivir to be held on December 06, 2025 in the city of Jaipur, Rajasthan"
$ws.Range("G4").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1764675345527.pdf"
$ws.Range("H4").Value = "1764675345527.pdf"
$ws.Range("I4").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Press Release/2025/December/1764675345527.pdf"

# Row 5: Consultation Paper on Review of Master Circular for Forei...
$ws.Range("A5").Value = "SEBI"
$ws.Range("B5").Value = "Consultation Paper"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "2025"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "December"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2025-12-05"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "Consultation Paper on Review of Master Circular for Foreign Portfolio Investors (FPIs) and Designated Depository Participants (DDPs)"
$ws.Range("G5").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1764943412112.pdf"
$ws.Range("H5").Value = "1764943412112.pdf"
$ws.Range("I5").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Consultation Paper/2025/December/1764943412112.pdf"

# Row 6: Consultation Paper on Review of existing position limits ...
$ws.Range("A6").Value = "SEBI"
$ws.Range("B6").Value = "Consultation Paper"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "2025"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "December"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2025-12-04"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "Consultation Paper on Review of existing position limits for Trading Members in Equity Derivatives Segment"
$ws.Range("G6").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1764859538149.pdf"
$ws.Range("H6").Value = "1764859538149.pdf"
$ws.Range("I6").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Consultation Paper/2025/December/1764859538149.pdf"

# Row 7: Securities and Exchange Board of India (Intermediaries​) ...
$ws.Range("A7").Value = "SEBI"
$ws.Range("B7").Value = "Regulations"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "2025"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "December"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2025-12-05"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "Securities and Exchange Board of India (Intermediaries​) (Third Amendment) Regulations, 2025"
$ws.Range("G7").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765172737392.pdf"
$ws.Range("H7").Value = "1765172737392.pdf"
$ws.Range("I7").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Regulations/2025/December/1765172737392.pdf"

# Row 8: Securities and Exchange Board of India (Substantial Acqui...
$ws.Range("A8").Value = "SEBI"
$ws.Range("B8").Value = "Regulations"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "2025"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "December"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2025-12-05"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "Securities and Exchange Board of India (Substantial Acquisition of Shares and Takeovers) (Amendment) Regulations, 2025"
$ws.Range("G8").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765172337309.pdf"
$ws.Range("H8").Value = "1765172337309.pdf"
$ws.Range("I8").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Regulations/2025/December/1765172337309.pdf"

# Row 9: Securities and Exchange Board of India (Merchant Bankers)...
$ws.Range("A9").Value = "SEBI"
$ws.Range("B9").Value = "Regulations"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "2025"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "December"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2025-12-05"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "Securities and Exchange Board of India (Merchant Bankers) (Amendment) Regulations, 2025"
$ws.Range("G9").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765171459782.pdf"
$ws.Range("H9").Value = "1765171459782.pdf"
$ws.Range("I9").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Regulations/2025/December/1765171459782.pdf"

# Row 10: Securities and Exchange Board of India (Share Based Emplo...
$ws.Range("A10").Value = "SEBI"
$ws.Range("B10").Value = "Regulations"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "2025"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "December"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2025-12-04"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = "Securities and Exchange Board of India (Share Based Employee Benefits and Sweat Equity) (Second Amendment) Regulations, 2025"
$ws.Range("G10").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765172885621.pdf"
$ws.Range("H10").Value = "1765172885621.pdf"
$ws.Range("I10").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Regulations/2025/December/1765172885621.pdf"

# Row 11: Securities and Exchange Board of India (Foreign Venture C...
$ws.Range("A11").Value = "SEBI"
$ws.Range("B11").Value = "Regulations"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "2025"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "December"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2025-12-03"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = "Securities and Exchange Board of India (Foreign Venture Capital Investors) (Amendment) Regulations, 2025"
$ws.Range("G11").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1764762979981.pdf"
$ws.Range("H11").Value = "1764762979981.pdf"
$ws.Range("I11").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Regulations/2025/December/1764762979981.pdf"

# Row 12: Securities and Exchange Board of India (Foreign Portfolio...
$ws.Range("A12").Value = "SEBI"
$ws.Range("B12").Value = "Regulations"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "2025"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "December"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2025-12-03"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = "Securities and Exchange Board of India (Foreign Portfolio Investors) (Second Amendment) Regulations, 2025"
$ws.Range("G12").Value = "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1764762906502.pdf"
$ws.Range("H12").Value = "1764762906502.pdf"
$ws.Range("I12").Value = "/Users/admin/Downloads/Tejomaya_pdfs_test/Akshayam Data/SEBI/Regulations/2025/December/1764762906502.pdf"

# Re-add hyperlinks for column G (PDF_URL) rows 2-12, restoring the Hyperlink style
$ws.Hyperlinks.Add($ws.Cells.Item(2,7), "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765190360681.pdf")
$ws.Range("G2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(3,7), "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1764758072366.pdf")
$ws.Range("G3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(4,7), "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1764675345527.pdf")
$ws.Range("G4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(5,7), "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1764943412112.pdf")
$ws.Range("G5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(6,7), "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1764859538149.pdf")
$ws.Range("G6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(7,7), "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765172737392.pdf")
$ws.Range("G7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(8,7), "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765172337309.pdf")
$ws.Range("G8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(9,7), "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765171459782.pdf")
$ws.Range("G9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(10,7), "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1765172885621.pdf")
$ws.Range("G10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(11,7), "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1764762979981.pdf")
$ws.Range("G11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(12,7), "https://www.sebi.gov.in/sebi_data/attachdocs/dec-2025/1764762906502.pdf")
$ws.Range("G12").Style = "Hyperlink"
